$wb = $excel.ActiveWorkbook

# Sheet "展览" (F column updates)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 592
$ws1.Range("F10").Value = 402
$ws1.Range("F11").Value = 360
$ws1.Range("F12").Value = 789
$ws1.Range("F13").Value = 787
$ws1.Range("F16").Value = 1547
$ws1.Range("F17").Value = 1547
$ws1.Range("F18").Value = 1092
$ws1.Range("F20").Value = 1367
$ws1.Range("F21").Value = 167
$ws1.Range("F22").Value = 367
$ws1.Range("F25").Value = 117
$ws1.Range("F26").Value = 6764
$ws1.Range("F27").Value = 5186
$ws1.Range("F28").Value = 10
$ws1.Range("F29").Value = 152
$ws1.Range("F32").Value = 218
$ws1.Range("F35").Value = 31
$ws1.Range("F37").Value = 1322
$ws1.Range("F38").Value = 201
$ws1.Range("F39").Value = 261
$ws1.Range("F40").Value = 635
$ws1.Range("F43").Value = 271
$ws1.Range("F44").Value = 154
$ws1.Range("F48").Value = 106

# Sheet "演出" (F column updates)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 15

# Sheet "本地生活" (F column updates)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2481
$ws3.Range("F4").Value = 220

# Sheet "全部类型" (F column updates)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 220
$ws4.Range("F9").Value = 592
$ws4.Range("F13").Value = 15
$ws4.Range("F14").Value = 402
$ws4.Range("F15").Value = 360
$ws4.Range("F16").Value = 789
$ws4.Range("F17").Value = 787
$ws4.Range("F20").Value = 1547
$ws4.Range("F21").Value = 1547
$ws4.Range("F22").Value = 1092
$ws4.Range("F24").Value = 367
$ws4.Range("F26").Value = 117
$ws4.Range("F29").Value = 6764
$ws4.Range("F30").Value = 5186
$ws4.Range("F31").Value = 218
$ws4.Range("F32").Value = 31
$ws4.Range("F33").Value = 1322
$ws4.Range("F34").Value = 201
$ws4.Range("F36").Value = 261
$ws4.Range("F38").Value = 635
$ws4.Range("F43").Value = 271
$ws4.Range("F47").Value = 106
